{"js": "// Update the date paragraph and every division-problem cell in the table.\n// The document is a single title paragraph (\"YYYY-MM-DD Weekday\") followed\n// by one table containing 5 \"problem\" rows (each with 5 cells) interleaved\n// with 15 blank rows. We update the title text and the 25 problem cells\n// in document order, matching the authored diff positionally.\n\nconst body = context.document.body;\n\n// 1) Title paragraph: \"2025-01-22 Wednesday\" -> \"2025-01-23 Thursday\"\nconst titleParas = body.paragraphs;\ntitleParas.load(\"items\");\nawait context.sync();\n\nif (titleParas.items.length > 0) {\n  const titlePara = titleParas.items[0];\n  titlePara.load(\"text\");\n  await context.sync();\n  if (titlePara.text.indexOf(\"2025-01-22 Wednesday\") !== -1) {\n    titlePara.insertText(\"2025-01-23 Thursday\", \"Replace\");\n  }\n}\n\n// 2) Table cells: replace the 25 problem values, in row-major order, with\n// their new values (row indices with content: 0, 4, 8, 12, 16).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValuesByRow = {\n  0: [\"23\u00f76=\", \"23\u00f75=\", \"63\u00f75=\", \"24\u00f75=\", \"53\u00f75=\"],\n  4: [\"32\u00f74=\", \"73\u00f75=\", \"62\u00f77=\", \"73\u00f72=\", \"91\u00f75=\"],\n  8: [\"93\u00f75=\", \"18\u00f75=\", \"29\u00f74=\", \"58\u00f77=\", \"37\u00f73=\"],\n  12: [\"36\u00f79=\", \"73\u00f73=\", \"96\u00f77=\", \"23\u00f72=\", \"43\u00f79=\"],\n  16: [\"31\u00f73=\", \"92\u00f76=\", \"74\u00f72=\", \"62\u00f78=\", \"59\u00f72=\"],\n};\n\nfor (const rowIndexStr of Object.keys(newValuesByRow)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const values = newValuesByRow[rowIndex];\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date paragraph and every division-problem cell in the table.\n# The document is a single title paragraph (\"YYYY-MM-DD Weekday\") followed\n# by one table containing 5 \"problem\" rows (each with 5 cells) interleaved\n# with 15 blank rows. We update the title text and the 25 problem cells,\n# matching the authored diff positionally.\n\n$d = $word.ActiveDocument\n\n# 1) Title paragraph: \"2025-01-22 Wednesday\" -> \"2025-01-23 Thursday\"\n$find = $d.Content.Find\n$find.Text = \"2025-01-22 Wednesday\"\n$find.Replacement.Text = \"2025-01-23 Thursday\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2) Table cells: replace the 25 problem values, in row-major order, with\n# their new values (1-indexed content rows: 1, 5, 9, 13, 17).\n$t = $d.Tables.Item(1)\n\n$t.Cell(1,1).Range.Text = \"23\u00f76=\"\n$t.Cell(1,2).Range.Text = \"23\u00f75=\"\n$t.Cell(1,3).Range.Text = \"63\u00f75=\"\n$t.Cell(1,4).Range.Text = \"24\u00f75=\"\n$t.Cell(1,5).Range.Text = \"53\u00f75=\"\n\n$t.Cell(5,1).Range.Text = \"32\u00f74=\"\n$t.Cell(5,2).Range.Text = \"73\u00f75=\"\n$t.Cell(5,3).Range.Text = \"62\u00f77=\"\n$t.Cell(5,4).Range.Text = \"73\u00f72=\"\n$t.Cell(5,5).Range.Text = \"91\u00f75=\"\n\n$t.Cell(9,1).Range.Text = \"93\u00f75=\"\n$t.Cell(9,2).Range.Text = \"18\u00f75=\"\n$t.Cell(9,3).Range.Text = \"29\u00f74=\"\n$t.Cell(9,4).Range.Text = \"58\u00f77=\"\n$t.Cell(9,5).Range.Text = \"37\u00f73=\"\n\n$t.Cell(13,1).Range.Text = \"36\u00f79=\"\n$t.Cell(13,2).Range.Text = \"73\u00f73=\"\n$t.Cell(13,3).Range.Text = \"96\u00f77=\"\n$t.Cell(13,4).Range.Text = \"23\u00f72=\"\n$t.Cell(13,5).Range.Text = \"43\u00f79=\"\n\n$t.Cell(17,1).Range.Text = \"31\u00f73=\"\n$t.Cell(17,2).Range.Text = \"92\u00f76=\"\n$t.Cell(17,3).Range.Text = \"74\u00f72=\"\n$t.Cell(17,4).Range.Text = \"62\u00f78=\"\n$t.Cell(17,5).Range.Text = \"59\u00f72=\"\n"}
